$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time value for the remaining row
$ws.Range("C2").Value = "15:49:41"

# Remove the second data row (vipin / 15:21:46) entirely
$ws.Rows(3).Delete()
